$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Range("A5").Value = 42604.891377314816

$ws.Range("B5").Value = "Bag"
$ws.Range("C5").Value = 10073
$ws.Range("D5").Value = 12103
$ws.Range("E5").Value = 1369
$ws.Range("F5").Value = 191
$ws.Range("G5").Value = 141
$ws.Range("H5").Value = 57
$ws.Range("I5").Value = 42
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 20
$ws.Range("L5").Value = 16
$ws.Range("M5").Value = 83
